# Auto-update from RStudio [gert]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap Story_title/Link between row 3 and row 4 ---
$b3New = @'
Akron women arrested in major meth bust in New York - FOX 8 News
'@
$e3New = @'
https://tinyurl.com/26ghp4vs
'@
$b4New = @'
Two Charged Following Three-Year Drug Investigation - Post Journal
'@
$e4New = @'
https://tinyurl.com/27k256qt
'@
$ws.Range("B3").Value = $b3New
$ws.Range("E3").Value = $e3New
$ws.Range("B4").Value = $b4New
$ws.Range("E4").Value = $e4New

# --- Update row 17 (Story_title/Link) ---
$b17New = @'
Rolison Sentenced For Having Meth At Jail - Post Journal
'@
$e17New = @'
https://tinyurl.com/2amvc7b2
'@
$ws.Range("B17").Value = $b17New
$ws.Range("E17").Value = $e17New

# --- Row 18 takes on what used to be row 19's content ---
$a18New = @'
WJTN News Headlines for Tues., Apr. 1, 2025
'@
$b18New = @'
2 keyword references (WJTN)
'@
$c18New = @'
Jamestown man convicted of rape and criminal sex act resentenced in County Court... A city man convicted of first-degree rape five years ago has now been resentenced in Chautauqua County Court to 18 years in prison and 25-years post-release supervision.  District Attorney Jason Schmidt says County Court Judge David Foley Monday handed down the sentence to 40 year-old Dennis Park.  Schmidt says Park also received an 18 year sentence... and, 25 years post-release supervision for his conviction for first-degree criminal sex act.  He says both sentences will be served concurrently.  Schmidt says Park was previously sentenced on the two charges in March 2020... but, the state Appellate Division ruled that the sentence was illegal, and ordered it to be returned to Chautauqua County Court in November of last year.  Schmidt says Park had the option to be resentenced or withdraw his plea... and, he chose to be resentenced.  The original incident occurred in August 2018 and the investigation was handled by Jamestown Police.Rolison to spend another 2-to-4 years in State Prison for possessing prison contraband... The Jamestown man sentenced to up to 40 years in prison for the motor vehicle deaths of a city teen... and, a woman in the town of Arkwright will spend another 2 to 4 years for promoting prison contraband.  Chautauqua County District Attorney Jason Schmidt has announced that 61 year-old Randall Rolison was sentenced Monday on his conviction for first-degree promoting prison contraband -- a Felony -- to spend the additional term consecutively to the lengthy prison sentences handed down just over a year ago in County Court.  Schmidt says Rolison was sentenced to the combined prison term for his second-degre manslaughter conviction for the hit and run death of 15-year-old Lexi Hughan.  He also was sentenced for aggravated vehicular homicide for the December 2022 car accident when his vehicle broadsided a pick-up truck driven by Gary Kraemer that killed his wife, Linda, and severely injured him.  In both incidents... the D.A. says Rolison was under the influence of methamphetamine.  In March of last year... Rolison was sentenced in county court... but, before he was transferred to a state facility... corrections officers at the County jail found him in possession of methamphetamine following a jail visit, leading to this most recent conviction.  “I can’t think of a more horrendous example of how illegal drugs like meth have destroyed innocent lives and caused needless and profound suffering for all of us here in the County.  It defies common sense to think that Mr. Rolison was on meth when he killed poor Lexi, chose to use meth again later that same year when he killed Mrs. Kraemer, and turned to meth again while in our jail before leaving for state prison."City Council rejects reallocating $640,000 in HOME funding for Gateway Lofts project...The Jamestown City Council has rejected the reallocation of $640,000 for remediation work at the site of the proposed Gateway Lofts housing project on Water Street.  However... it's believed that the estimated $77-million project will be able to move forward with money from confirmed investors.  The project was debated again at last night's Council voting session prior to lawmakers rejecting the move, 5-to-4.  Housing Committee Chairman Bill Reynolds says he'd rather see that money be used to improve existing housing... because since the investors have a lot of money for upfront costs, they don't need that funding.  Councilwoman Regina Brackman pointed out that City Development Director Crystal Surdyk will have to move quickly on finding somewhere else to use that HOME Program money... but, she also hoped to see it used to improve the city's housing stock.  Council President Tony Dolce voted for the reallocation and expressed similar feelings, in hoping the money can be used for existing housing improvements.  The city's Department of Development had packaged the money from the past five year's remaining HOME Program funding to support the project.  Surdyk says it will be difficult to find another use for that funding because the federal government puts very tight restrictions on it... and, right now... they don't have a project for that amount of money.Debate fairly split among those favoring Gateway Lofts funding vs. those against it...Prior to last night's vote on the funding for the Gateway Lofts project... there were several comments made during a public hearing... and, more during the later voting session.  Former City Councilman Tom Nelson... who now represents that area on the county legislature... says he was also on the city's Planning Commission when the project was first proposed several years ago.  He urged it's passage... says he's seen the project "really evolve from what was an okay project to one a better, and better project."  Nelson feels it's something the city "needs."  A number of people spoke against the project... and, one... Jamestown Police Lieutenant Sam Piazza... felt city police would be overwhelmed by the number of calls there and elsewhere in the city.  Piazza... who is also president of the policeman's union... says JPD doesn't have enough officers as it is.  The man heading up the project for Southern Tier Environment's for Living... Steve Ald... reiterated that the facility would have "24/7" security with a front desk staffed 24-hours each day with four-and-a-half people.  However... Piazza said -- if their security is unable to handle it... JPD would be called. 
'@
$e18New = @'
https://mediaonegroupradio.com/news-and-closings/local-news-headlines/786700/wjtn-news-headlines-for-tues-apr-1-2025
'@
$f18New = @'
...nbsp;In both incidents... the D.A. says Rolison was under the influence of **methamphetamine**.  In March of last year... Rolison was sentenced in county court... b...

...lity... corrections officers at the County jail found him in possession of **methamphetamine** following a jail visit, leading to this most recent conviction.  &ldq...
'@
$g18New = @'
News_search_MediaOne_WJTN.R
'@
$ws.Range("A18").Value = $a18New
$ws.Range("B18").Value = $b18New
$ws.Range("C18").Value = $c18New
$ws.Range("E18").Value = $e18New
$ws.Range("F18").Value = $f18New
$ws.Range("G18").Value = $g18New

# --- Delete row 19 entirely (its content was merged into row 18) ---
$ws.Rows(19).Delete()
